# Edit summary:
#  1. Three tables (on slides 14, 15 and 16) get their table style switched
#     from {89D5681D-51B3-4CD6-9506-B016307E953B} to
#     {814365C5-A3B9-49BD-8AA8-394D83280B63}.
#  2. The presentation's design theme ("Integral" / Red Violet color scheme,
#     stored in ppt/theme/theme2.xml - the theme actually used by the slide
#     master) is swapped back to the plain "Office Theme" color values that
#     previously lived in ppt/theme/theme1.xml (used only by the notes
#     master). I.e. the 12 theme colors get swapped between the two
#     palettes.

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables ------------------------------------------------
$newTableStyle = "{814365C5-A3B9-49BD-8AA8-394D83280B63}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2. Swap the active theme's color scheme to the "Office Theme" palette ------
$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

# Order used by ThemeColorScheme.Item(n): dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
# Values below are the "Office Theme" palette (previously ppt/theme/theme1.xml),
# expressed as COM RGB() (0x00BBGGRR) integers so Item(n).RGB = value round-trips
# to the matching srgbClr hex.
$officeThemeRgb = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeRgb[$i - 1]
}
